# "clean up namings across simulations"
# On the "Sim2" sheet, a new condition row (equal split, n=160 per arm,
# avg(10x0.75)=0.75) is inserted above the existing condition rows, and a
# typo in one of the existing condition labels (...0.76 -> ...0.75) is
# fixed at the same time.

$wb = $excel.ActiveWorkbook

# ---- Sim2 sheet -----------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sim2")
$ws2.Activate()

# Insert a new row above row 2, shifting the existing condition rows
# (and every row below them) down by one. Use the row that ends up
# below (old row 2, now row 3) as the formatting template for the new
# row.
$ws2.Rows("2:2").Insert(-4121)
$ws2.Range("A3:M3").Copy()
$ws2.Range("A2:M2").PasteSpecial(-4122)

# Populate the newly inserted row with the new condition.
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "avg(10x0.75)=0.75"
$ws2.Range("C2:G2").Value = 160
$ws2.Range("H2:L2").Value = 160
$ws2.Range("M2").Value = "mean=median"

# The row that used to be row 2 (now row 3) keeps its own data, but its
# row height changes from the old 19.8 to the common 15.
$ws2.Rows("3:3").RowHeight = 15

# Fix the typo in the condition that is now on row 5
# ("...0.76" -> "...0.75").
$ws2.Range("B5").Value = "avg(5x0.60, 5x0.90)=0.75"

# The row that used to be the blank row 6 now carries the condition
# that used to be on row 5; renumber its index column.
$ws2.Range("A6").Value = 5

$ws2.Range("E13").Select()

# ---- Sim1 sheet -------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sim1")
$ws1.Activate()
$ws1.Range("C6:L6").Select()

# Leave "Sim2" as the active sheet/tab, matching the saved state.
$ws2.Activate()
